# Correct the plebdex BTC ticker: "BTCUSD" -> "BTC-USD"
# The ticker lives in the "2023" worksheet and is referenced (via the
# shared string table) by the other yearly sheets as well, so simply
# updating the shared text here is enough for Excel to keep every sheet
# in sync.
$wb = $excel.ActiveWorkbook

$ws2023 = $wb.Worksheets.Item("2023")

# Locate the BTCUSD ticker cell within column A and fix its spelling.
$tickerRange = $ws2023.Range("A1:A10")
$btcCell = $tickerRange.Find("BTCUSD")
if ($btcCell -ne $null) {
    $btcCell.Value = "BTC-USD"
} else {
    $ws2023.Range("A7").Value = "BTC-USD"
}

# Make the "2023" tab the active sheet/selection, matching the saved view
# state captured in the workbook.
$ws2023.Activate()
$ws2023.Range("A7").Select()

$wb.Save()
